$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AY2").Value = 11975530.496
$ws.Range("AZ2").Value = 12455377.92
$ws.Range("AY3").Value = 5423211.008
$ws.Range("AZ3").Value = 5419362.816
$ws.Range("AY4").Value = 1252747.008
$ws.Range("AZ4").Value = 984729.9840000001
$ws.Range("AY5").Value = 607524.992
$ws.Range("AZ5").Value = 415880.992
$ws.Range("AY6").Value = 335060.992
$ws.Range("AZ6").Value = 348539.008
$ws.Range("AY7").Value = 91085
$ws.Range("AZ7").Value = 109429
$ws.Range("AY8").Value = 1267
$ws.Range("AZ8").Value = 894
$ws.Range("AY9").Value = 2765644.032
$ws.Range("AZ9").Value = 3160880.128
$ws.Range("AY11").Value = 256248.992
$ws.Range("AZ11").Value = 337132
$ws.Range("AY12").Value = 113633
$ws.Range("AZ12").Value = 61878
$ws.Range("AY13").Value = 5761196.032
$ws.Range("AZ13").Value = 6238730.752
$ws.Range("AY14").Value = 16828
$ws.Range("AZ14").Value = 118984
$ws.Range("AY15").Value = 586497.9840000001
$ws.Range("AZ15").Value = 599286.976
$ws.Range("AY16").Value = 0
$ws.Range("AZ16").Value = 0
$ws.Range("AY18").Value = 4451140.096
$ws.Range("AZ18").Value = 4733943.808
$ws.Range("AY20").Value = 682091.008
$ws.Range("AZ20").Value = 766270.0159999999
$ws.Range("AY21").Value = 24639
$ws.Range("AZ21").Value = 20246
$ws.Range("AY22").Value = 791123.008
$ws.Range("AZ22").Value = 797283.968
$ws.Range("AY23").Value = 549948.032
$ws.Range("AZ23").Value = 564406.976
$ws.Range("AY24").Value = 159008
$ws.Range("AZ24").Value = 151156.992
$ws.Range("AY26").Value = 82167
$ws.Range("AZ26").Value = 81720
$ws.Range("AY27").Value = 0
$ws.Range("AZ27").Value = 0
$ws.Range("AY28").Value = 11975530.496
$ws.Range("AZ28").Value = 12455377.92
$ws.Range("AY29").Value = 4253090.048
$ws.Range("AZ29").Value = 4499202.048
$ws.Range("AY30").Value = 3393500.928
$ws.Range("AZ30").Value = 3486223.872
$ws.Range("AY31").Value = 90187
$ws.Range("AZ31").Value = 131958
$ws.Range("AY32").Value = 70244
$ws.Range("AZ32").Value = 70004
$ws.Range("AY33").Value = 87934
$ws.Range("AZ33").Value = 124225
$ws.Range("AY34").Value = 3854
$ws.Range("AZ34").Value = 15670
$ws.Range("AY37").Value = 0
$ws.Range("AZ37").Value = 0
$ws.Range("AY38").Value = 607369.9840000001
$ws.Range("AZ38").Value = 671121.024
$ws.Range("AY39").Value = 6609331.2
$ws.Range("AZ39").Value = 6831449.088
$ws.Range("AY40").Value = 5800022.016
$ws.Range("AZ40").Value = 5891651.072
$ws.Range("AY42").Value = 273
$ws.Range("AY43").Value = 0
$ws.Range("AZ43").Value = 0
$ws.Range("AY47").Value = 0
$ws.Range("AZ47").Value = 0
$ws.Range("AY48").Value = 809036.032
$ws.Range("AZ48").Value = 939798.0159999999
$ws.Range("AY49").Value = 273
$ws.Range("AZ49").Value = 236
$ws.Range("AY51").Value = 1112835.968
$ws.Range("AZ51").Value = 1124491.008
$ws.Range("AY52").Value = 492708
$ws.Range("AZ52").Value = 597539.968
$ws.Range("AY53").Value = 43375
$ws.Range("AZ53").Value = 43375
$ws.Range("AY54").Value = 111
$ws.Range("AZ54").Value = 108
$ws.Range("AY55").Value = 541620.992
$ws.Range("AZ55").Value = 478081.984
$ws.Range("AY56").Value = -5109
$ws.Range("AZ56").Value = 11000
$ws.Range("AY57").Value = 40130
$ws.Range("AZ57").Value = 0
$ws.Range("AY58").Value = 655587.008
$ws.Range("AZ58").Value = 809446.976
$ws.Range("AY59").Value = -279056
$ws.Range("AZ59").Value = -327936
$ws.Range("AY60").Value = 376531.04
$ws.Range("AZ60").Value = 481511.008
$ws.Range("AY61").Value = -445104.064
$ws.Range("AZ61").Value = -417344
$ws.Range("AY62").Value = 63027
$ws.Range("AZ62").Value = 67597
$ws.Range("AY63").Value = -122074
$ws.Range("AZ63").Value = -109736
$ws.Range("AY64").Value = -199014.976
$ws.Range("AZ64").Value = -215304.992
$ws.Range("AY65").Value = -34369.992
$ws.Range("AZ65").Value = -33708
$ws.Range("AY66").Value = 31900
$ws.Range("AZ66").Value = 39930
$ws.Range("AY67").Value = -217934.032
$ws.Range("AZ67").Value = -182832.992
$ws.Range("AY68").Value = 33362
$ws.Range("AZ68").Value = 16711
$ws.Range("AY69").Value = -68573
$ws.Range("AZ69").Value = 64167
$ws.Range("AY70").Value = -6840
$ws.Range("AZ70").Value = 780
$ws.Range("AY71").Value = 855
$ws.Range("AZ71").Value = 1085
$ws.Range("AY72").Value = -7695
$ws.Range("AZ72").Value = -305
$ws.Range("AY73").Value = -75412.992
$ws.Range("AZ73").Value = 64947
$ws.Range("AY74").Value = 125316
$ws.Range("AZ74").Value = -9881
$ws.Range("AY76").Value = -6083
$ws.Range("AZ76").Value = -7949
$ws.Range("AY77").Value = 0
$ws.Range("AZ77").Value = 0
$ws.Range("AY79").Value = 43820.008
$ws.Range("AZ79").Value = 47117
